$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=84905; B="André da Cruz";          C="Engenharia";        D="Doença";              E=6; F=45078; G=5534.56},
    @{Row=3;  A=9135;  B="Heitor Vieira";           C="Jurídico";          D="Doença";              E=8; F=45097; G=12377.19},
    @{Row=4;  A=98463; B="Dr. Enrico Silveira";     C="Marketing";         D="Viagem de negócios";  E=6; F=45093; G=4311.16},
    @{Row=5;  A=68518; B="Sarah Silveira";          C="Jurídico";          D="Outros";              E=1; F=45100; G=2740.67},
    @{Row=6;  A=87929; B="Sophia Ramos";            C="Operações";         D="Viagem de negócios";  E=4; F=45105; G=9819.39},
    @{Row=7;  A=41461; B="Maria Vitória Silveira";  C="Financeiro";        D="Consulta médica";     E=8; F=45090; G=2873.79},
    @{Row=8;  A=49469; B="Kamilly Cardoso";         C="Engenharia";        D="Consulta médica";     E=2; F=45094; G=8321.04},
    @{Row=9;  A=8134;  B="Felipe Lopes";            C="Vendas";            D="Doença";              E=6; F=45100; G=7893.41},
    @{Row=10; A=83606; B="Henrique Rezende";        C="Marketing";         D="Consulta médica";     E=5; F=45085; G=9065.17},
    @{Row=11; A=39133; B="Alice Fogaça";            C="Recursos Humanos";  D="Consulta médica";     E=1; F=45082; G=9581.56}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
